# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets (row 3 -> the
# fa496268-6053-49df-8127-4f8d4c5800aa.* entry), as part of regenerating
# the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-28 05:30:45"
$wsZhCn.Range("G3").Value = "2016-01-28 05:31:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-28 05:30:57"
$wsDeDe.Range("G3").Value = "2016-01-28 05:32:06"
